$d = $word.ActiveDocument

# Locate the "Instructor" run (the heading before the "Instructor:" label)
$rng = $d.Content.Duplicate
$rng.Find.Execute("Instructor", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Collapse to the end of the found text and insert a new run containing "s"
$rng.Collapse(0)
$rng.InsertAfter("s")
$rng.Font.Name = "Century Gothic"
$rng.Font.NameFarEast = "Century Gothic"
$rng.Font.NameAscii = "Century Gothic"
$rng.Font.NameOther = "Century Gothic"
$rng.Font.Bold = $true
$rng.Font.Size = 10
